$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.488.89'
$ws.Range('E2').Value = '  -5.65%  '
$ws.Range('D3').Value = '3.267.73'
$ws.Range('E3').Value = '  -6.03%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '560.74'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -4.02%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '126.53'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.74%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').Value = '3.265.81'
$ws.Range('E8').Value = '  -6.06%  '
$ws.Range('E9').Value = '  -1.89%  '
$ws.Range('E10').Value = '  -4.35%  '
$ws.Range('E11').Value = '  -4.87%  '
$ws.Range('E12').Value = '  -4.06%  '
$ws.Range('D13').Value = '3.829.91'
$ws.Range('E13').Value = '  -5.92%  '
$ws.Range('E14').Value = '  -0.13%  '
$ws.Range('D15').Value = '3.277.83'
$ws.Range('E15').Value = '  -5.76%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000166'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -6.31%  '
$ws.Range('D17').Value = '60.669.51'
$ws.Range('E17').Value = '  -5.39%  '
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.61'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.39%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.21'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.44%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '8.95'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -10.33%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '349.88'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -8.96%  '
$ws.Range('E23').Value = '  -3.36%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.00'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').Value = '3.401.95'
$ws.Range('E25').Value = '  -5.96%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '69.01'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -7.74%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000106'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -5.02%  '
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.17'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.76%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.41'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.93%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.77'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.07%  '
$ws.Range('E32').Value = '  -6.31%  '
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('E34').Value = '  -1.73%  '
$ws.Range('D35').Value = '3.299.65'
$ws.Range('E35').Value = '  -5.91%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '22.54'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.76%  '
$ws.Range('E37').Value = '  +0.49%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.75'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.47'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -2.01%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '158.86'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.33%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0748'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.41%  '
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '40.89'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.18%  '
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.735'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -7.82%  '
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.54'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.74%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '22.22'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -4.64%  '
$ws.Range('E49').Value = '  -0.81%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.857'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -4.84%  '
$ws.Range('E51').Value = '  +3.80%  '
